# "time update, code bug fix and CSV map generation"
#
# - Rename Sheet2 -> subsequences and populate it with a point -> A/B/C
#   lookup table (same layout/formulas as the "sequences" sheet's header +
#   first 4 data rows).
# - Update the saved selection on "sequences" to A1:E5 (no longer the
#   active tab).
# - Make "subsequences" the active tab, with cell G13 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- rename Sheet2 -> subsequences -------------------------------------
$ws2.Name = "subsequences"

# --- center-align the table region, like the source sheet --------------
$ws2.Range("A1:E5").HorizontalAlignment = -4108   # xlCenter

# --- header row ----------------------------------------------------------
$ws2.Range("A1").Value = "Point"
$ws2.Range("B1").Value = "Label"
$ws2.Range("C1").Value = "A"
$ws2.Range("D1").Value = "B"
$ws2.Range("E1").Value = "C"

# --- data rows (mirrors "sequences" rows 2-5) ---------------------------
$ws2.Range("A2").Value = 1
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("B2").Formula = "=IF(C2=1,""A"",IF(D2=1,""B"",IF(E2=1,""C"",""X"")))"

$ws2.Range("A3").Value = 2
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0
$ws2.Range("B3").Formula = "=IF(C3=1,""A"",IF(D3=1,""B"",IF(E3=1,""C"",""X"")))"

$ws2.Range("A4").Value = 3
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 1
$ws2.Range("B4").Formula = "=IF(C4=1,""A"",IF(D4=1,""B"",IF(E4=1,""C"",""X"")))"

$ws2.Range("A5").Value = 4
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 1
$ws2.Range("B5").Formula = "=IF(C5=1,""A"",IF(D5=1,""B"",IF(E5=1,""C"",""X"")))"

# --- selections / active tab --------------------------------------------
# "sequences" keeps a saved selection of A1:E5 but is no longer the
# tab in front.
$ws1.Range("A1:E5").Select()

# "subsequences" becomes the active/front sheet, with G13 selected.
$ws2.Activate()
$ws2.Range("G13").Select()
